# Swap the data between row 4 and row 5 for the columns that differ:
# A (Id), B (Taxonsorteringsordning), E (TaxonId), F (Artnamn),
# G (Vetenskapligt namn), H (Auktor), Q (Ost), R (Nord)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $cellRow4 = $ws.Range($col + "4")
    $cellRow5 = $ws.Range($col + "5")

    $val4 = $cellRow4.Value2
    $val5 = $cellRow5.Value2

    $cellRow4.Value2 = $val5
    $cellRow5.Value2 = $val4
}
